# "Generate Report for Handback" - refresh the localization-status report
# after a handback sync: update status text, refresh handback timestamps,
# clear the stale "not latest" error now that things are in sync, and
# widen the Status / Error Detail columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value     = "Handed back: in sync with en-US"
$dede.Range("C2").Value     = "Handed back: in sync with en-US"

# --- zh-cn: refresh Latest Handback DateTime, clear the stale error detail ---
$zhcn.Range("K2").Value = "2016-09-06 00:55:34"
$zhcn.Range("P2").Value = ""

# --- de-de: refresh Latest Handback DateTime, clear the stale error detail ---
$dede.Range("K2").Value = "2016-09-06 00:55:41"
$dede.Range("P2").Value = ""

# --- Column width refresh (auto-sized wider to fit new Status / narrower
#     now that Error Detail is empty) ---
$overview.Columns.Item(5).ColumnWidth = 29.2   # column E
$overview.Columns.Item(6).ColumnWidth = 29.2   # column F

$zhcn.Columns.Item(3).ColumnWidth  = 29.2      # column C (Status)
$zhcn.Columns.Item(16).ColumnWidth = 12.8      # column P (Error Detail)

$dede.Columns.Item(3).ColumnWidth  = 29.2      # column C (Status)
$dede.Columns.Item(16).ColumnWidth = 12.8      # column P (Error Detail)
